$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '42.025.65'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '2.251.79'
$ws.Range('E3').Value = '  -0.67%  '
Set-TextValue ($ws.Range('D4')) '0.999'
$ws.Range('E4').Value = '  -0.15%  '
Set-TextValue ($ws.Range('D5')) '306.59'
$ws.Range('E5').Value = '  +0.25%  '
Set-TextValue ($ws.Range('D6')) '97.42'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  -1.20%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -0.41%  '
Set-TextValue ($ws.Range('D10')) '35.20'
$ws.Range('E10').Value = '  -1.20%  '
Set-TextValue ($ws.Range('D11')) '0.0813'
$ws.Range('E11').Value = '  +2.28%  '
$ws.Range('E12').Value = '  +1.42%  '
Set-TextValue ($ws.Range('D13')) '6.73'
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('D14').Value = '2.599.31'
$ws.Range('E14').Value = '  -0.95%  '
Set-TextValue ($ws.Range('D15')) '14.44'
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').Value = '2.196.52'
$ws.Range('E16').Value = '  -3.13%  '
Set-TextValue ($ws.Range('D17')) '0.778'
$ws.Range('E17').Value = '  -1.91%  '
$ws.Range('D18').Value = '41.982.17'
$ws.Range('E18').Value = '  -0.50%  '
Set-TextValue ($ws.Range('D19')) '12.16'
$ws.Range('E19').Value = '  -2.72%  '
$ws.Range('D20').Value = '0.0₃0902'
$ws.Range('E20').Value = '  -0.74%  '
Set-TextValue ($ws.Range('D21')) '5.94'
$ws.Range('E21').Value = '  -0.36%  '
Set-TextValue ($ws.Range('D22')) '67.03'
$ws.Range('E22').Value = '  -0.67%  '
Set-TextValue ($ws.Range('D23')) '235.81'
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('E25').Value = '  -1.06%  '
Set-TextValue ($ws.Range('D26')) '0.995'
$ws.Range('E26').Value = '  -0.31%  '
Set-TextValue ($ws.Range('D27')) '37.93'
$ws.Range('E27').Value = '  +1.40%  '
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue ($ws.Range('D29')) '9.49'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue ($ws.Range('D30')) '2.11'
$ws.Range('E30').Value = '  +0.35%  '
Set-TextValue ($ws.Range('D31')) '166.61'
$ws.Range('E31').Value = '  +3.99%  '
Set-TextValue ($ws.Range('D32')) '0.999'
$ws.Range('E32').Value = '  -0.04%  '
Set-TextValue ($ws.Range('D33')) '5.17'
$ws.Range('E33').Value = '  -1.53%  '
Set-TextValue ($ws.Range('D34')) '17.51'
$ws.Range('E34').Value = '  +2.50%  '
$ws.Range('E35').Value = '  -3.52%  '
Set-TextValue ($ws.Range('D36')) '0.0720'
$ws.Range('E36').Value = '  -2.83%  '
$ws.Range('E37').Value = '  +1.24%  '
Set-TextValue ($ws.Range('D38')) '0.114'
$ws.Range('E38').Value = '  +0.11%  '
Set-TextValue ($ws.Range('D39')) '0.102'
$ws.Range('E39').Value = '  -2.58%  '
Set-TextValue ($ws.Range('D40')) '1.78'
$ws.Range('E40').Value = '  -2.24%  '
Set-TextValue ($ws.Range('D41')) '4.09'
$ws.Range('E41').Value = '  +0.83%  '
$ws.Range('D42').Value = '1.935.72'
$ws.Range('E42').Value = '  -2.88%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue ($ws.Range('D43')) '2.22'
$ws.Range('E43').Value = '  -8.89%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue ($ws.Range('D44')) '0.0281'
$ws.Range('E44').Value = '  -2.64%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue ($ws.Range('D45')) '18.57'
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('E47').Value = '  -3.06%  '
Set-TextValue ($ws.Range('D48')) '53.87'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('D49').Value = '2.471.71'
$ws.Range('E49').Value = '  -0.85%  '
Set-TextValue ($ws.Range('D50')) '71.20'
$ws.Range('E50').Value = '  -1.06%  '
Set-TextValue ($ws.Range('D51')) '91.08'
$ws.Range('E51').Value = '  -0.30%  '
